# Fixed naive component forecaster bug - Presentation state 11.02.
# This script reproduces a set of tiny floating point recalculations of the
# y_0_forecast (column C) / y_1_forecast (column E) vectors, plus clears a
# handful of early cells (C2, E2, C3, C4) that no longer have a forecast
# value available (insufficient history for those early vintages).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cells that are cleared entirely (no forecast value available) ---
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# --- Cells whose numeric forecast value is updated ---
$updates = @{
    "E3"  = 2.228397109637203
    "C5"  = -0.4626567965509643
    "E5"  = -0.2018858887078645
    "C6"  = -0.8792832172735632
    "E6"  = -0.1037449741509211
    "E7"  = 0.2108047537406454
    "C8"  = 0.9337833426867226
    "E9"  = 0.96842791562195
    "C10" = 2.791140000794257
    "E10" = 1.722110645261954
    "E11" = 1.698182372097512
    "C12" = 0.4451370000809529
    "E12" = 0.6480763427742176
    "C13" = -0.2674335569108788
    "E14" = 1.223618887196509
    "C15" = 2.038609866767915
    "E15" = 1.325176859452348
    "E17" = 1.522808462763692
    "E18" = 1.582979977679555
    "C19" = 2.246337373618967
    "E19" = 1.693557061600948
    "C20" = 2.2044495746113
    "E20" = 1.401113624217065
    "C21" = 1.777150434343544
    "E21" = 1.905564797014669
    "E22" = 2.357704431248386
    "C23" = 1.953801996162019
    "E23" = 2.2850393561338
    "C24" = 1.741137453897323
    "E24" = 2.08247707460909
    "E25" = 2.263031582094888
    "E26" = 2.36261304543155
    "C27" = 1.831523226563148
    "E27" = 2.069147625228918
    "C28" = 1.562095320687429
    "E28" = 1.845103901518907
    "E29" = 1.670042545405948
    "E30" = 0.7756897792100093
    "C31" = 0.9716673922242069
    "E31" = 1.221567067111828
    "E32" = 0.02570757229445331
    "C33" = -4.58200588280312
    "E33" = -5.994170454492931
    "C34" = -3.258619210312896
    "E34" = -1.049961713694159
    "C35" = -1.267247591471976
    "C37" = -0.4671761491813142
    "C38" = 0.4255262881966759
    "E38" = 2.734996705911397
    "C39" = 1.709116405894862
    "E39" = -0.3788276656591538
    "E42" = 1.813346177122321
    "C43" = 1.845950461732082
    "C44" = -0.2883789941992232
    "C45" = -0.1703190263132703
    "E45" = 1.429443202289193
    "C46" = -0.2814561130375703
    "E47" = 0.8057832540545151
    "C49" = -0.9450815092640896
    "E49" = 0.1338934428748884
    "C50" = -0.6470065423293758
    "E50" = 0.9227184786156251
    "C51" = 2.137819896031878
    "E51" = 0.3874019458836253
    "C52" = 1.431852292002245
    "C53" = 2.034296706251948
    "E53" = 0.6096213037243281
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
